$d = $word.ActiveDocument

# --- 1. Insert a new paragraph (with 3 runs) right after the existing
#        "cd ...GitHub\BuildingSync\schema\Java Files" paragraph and
#        before the bookmark ("_GoBack") paragraph. -----------------------

$cdPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "cd C:*BuildingSync\schema\Java Files*") {
        $cdPara = $p
        break
    }
}

$cdParaIndex = $cdPara.Index
$cdPara.Range.InsertParagraphAfter()
$newParaIndex = $cdParaIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$r = $newPara.Range
$r.Collapse(1)

$seg1 = "cd "
$seg2 = "C:\Users\Robert\Documents\NREL\BuildingSync\CTS\Java"
$seg3 = " Files"

# Type the first segment, split the paragraph (so the next segment
# lands in its own run), type the second segment, split again, then
# type the third segment. Finally merge the three temp paragraphs back
# into a single paragraph by deleting the two paragraph marks that were
# inserted - this keeps each segment in its own <w:r> run instead of
# letting same-formatted runs coalesce into one.
$r.InsertAfter($seg1)
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(1, 1) | Out-Null
$r.InsertAfter($seg2)
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(1, 1) | Out-Null
$r.InsertAfter($seg3)

$mergePara1 = $d.Paragraphs.Item($newParaIndex)
$mark1 = $d.Range($mergePara1.Range.End - 1, $mergePara1.Range.End)
$mark1.Delete()

$mergePara2 = $d.Paragraphs.Item($newParaIndex)
$mark2 = $d.Range($mergePara2.Range.End - 1, $mergePara2.Range.End)
$mark2.Delete()

# --- 2. Add one more empty paragraph among the trailing empty
#        paragraphs (3 -> 4), inserted right after the bookmark
#        ("_GoBack") paragraph, which immediately follows the paragraph
#        we just built. -------------------------------------------------

$bookmarkParaIndex = $newParaIndex + 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkPara.Range.InsertParagraphAfter()
